# Synchronized with the demo rep
# Applies the structural edit to the "RVL" sheet (sheet1):
#  - removes the MyMap range Param rows (fromRow/fromCol/toRow/toCol)
#  - renames the Functions/DfoXxx bootstrap actions to DFO/Xxx
#  - inserts a new "wait for grid / read cell / trim / output / tester / if" block
#  - inserts a new "End of If" block after the Save button click
#  - removes the trailing "Navigator/Close" action row
#
# All row numbers below refer to the ORIGINAL row numbers in the workbook.
# We work from the bottom of the sheet upward so that earlier row-number
# references used by later statements are never invalidated by row
# insertions/deletions that happen further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Remove the old "Action / Navigator / Close" row (original row 29)
# ---------------------------------------------------------------------
$ws.Range("A29").EntireRow.Delete()

# ---------------------------------------------------------------------
# 2) Insert the new "End of If" block right after the Save button click
#    (original row 25), i.e. before original row 26 (NavBarDashboard).
# ---------------------------------------------------------------------
$ws.Range("A26:A27").EntireRow.Insert()

$ws.Range("A26").Value = "End"
$ws.Range("B26").Value = "of If"
# A27 / row27 stays fully blank (spacer row)

# ---------------------------------------------------------------------
# 3) Insert the big new 14-row block after the blank spacer row
#    (original row 18), i.e. before original row 19 (SystemDefinedNewButton).
# ---------------------------------------------------------------------
$ws.Range("A19:A32").EntireRow.Insert()

# Row 19: Action / Global / DoWaitFor / objectId / objectid / custgroup_Grid
$ws.Range("B19").Value = "Action"
$ws.Range("C19").Value = "Global"
$ws.Range("D19").Value = "DoWaitFor"
$ws.Range("E19").Value = "objectId"
$ws.Range("F19").Value = "objectid"
$ws.Range("G19").Value = "custgroup_Grid"

# Row 20: Param / timeout / number / 120000
$ws.Range("B20").Value = "Param"
$ws.Range("E20").Value = "timeout"
$ws.Range("F20").Value = "number"
$ws.Range("G20").Value = "'120000"

# Row 21: Param / sleepTimeIfFound / number / 5000
$ws.Range("B21").Value = "Param"
$ws.Range("E21").Value = "sleepTimeIfFound"
$ws.Range("F21").Value = "number"
$ws.Range("G21").Value = "'5000"

# Row 22: Action / custgroup_Grid / GetCell / row / number / 0
$ws.Range("B22").Value = "Action"
$ws.Range("C22").Value = "custgroup_Grid"
$ws.Range("D22").Value = "GetCell"
$ws.Range("E22").Value = "row"
$ws.Range("F22").Value = "number"
$ws.Range("G22").Value = "'0"

# Row 23: Param / col / string / Customer group
$ws.Range("B23").Value = "Param"
$ws.Range("E23").Value = "col"
$ws.Range("F23").Value = "string"
$ws.Range("G23").Value = "Customer group"

# Row 24: Variable / CustGroup / variable / LastResult
$ws.Range("B24").Value = "Variable"
$ws.Range("E24").Value = "CustGroup"
$ws.Range("F24").Value = "variable"
$ws.Range("G24").Value = "LastResult"

# Row 25: Action / Global / DoTrim / str / variable / CustGroup
$ws.Range("A25").Value = ""
$ws.Range("B25").Value = "Action"
$ws.Range("C25").Value = "Global"
$ws.Range("D25").Value = "DoTrim"
$ws.Range("E25").Value = "str"
$ws.Range("F25").Value = "variable"
$ws.Range("G25").Value = "CustGroup"

# Row 26: Output / / variable / CustGroup
$ws.Range("A26").Value = ""
$ws.Range("B26").Value = "Output"
$ws.Range("F26").Value = "variable"
$ws.Range("G26").Value = "CustGroup"

# Row 27: blank spacer row

# Row 28: Action / Tester / Message / message / variable / CustGroup
$ws.Range("B28").Value = "Action"
$ws.Range("C28").Value = "Tester"
$ws.Range("D28").Value = "Message"
$ws.Range("E28").Value = "message"
$ws.Range("F28").Value = "variable"
$ws.Range("G28").Value = "CustGroup"

# Row 29: If / Param / param1 / variable / CustGroup
$ws.Range("A29").Value = "If"
$ws.Range("B29").Value = "Param"
$ws.Range("E29").Value = "param1"
$ws.Range("F29").Value = "variable"
$ws.Range("G29").Value = "CustGroup"

# Row 30: Condition / param1 != param2
$ws.Range("A30").Value = ""
$ws.Range("B30").Value = "Condition"
$ws.Range("D30").Value = "param1 != param2"

# Row 31: Param / param2 / Data / Name
$ws.Range("A31").Value = ""
$ws.Range("B31").Value = "Param"
$ws.Range("E31").Value = "param2"
$ws.Range("F31").Value = "Data"
$ws.Range("G31").Value = "Name"

# Row 32: # / If actions
$ws.Range("A32").Value = "#"
$ws.Range("B32").Value = "If actions"

# ---------------------------------------------------------------------
# 4) Rewrite the bootstrap "Functions" calls as "DFO" calls
#    (original rows 15, 16, 17 - unaffected by the inserts above, since
#    those all happened at row 19 or below).
# ---------------------------------------------------------------------
$ws.Range("C15").Value = "DFO"
$ws.Range("D15").Value = "Launch"

$ws.Range("C16").Value = "DFO"
$ws.Range("D16").Value = "PassWelcomeScreen"

$ws.Range("C17").Value = "DFO"
$ws.Range("D17").Value = "SearchPage"

# ---------------------------------------------------------------------
# 5) Remove the obsolete MyMap range Param rows (fromRow/fromCol/toRow/toCol)
#    (original rows 9-12).
# ---------------------------------------------------------------------
$ws.Range("A9:A12").EntireRow.Delete()
